$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "95÷2=47, 1"
$t.Cell(1, 2).Range.Text = "33÷6=5, 3"
$t.Cell(1, 3).Range.Text = "79÷8=9, 7"
$t.Cell(1, 4).Range.Text = "17÷8=2, 1"
$t.Cell(1, 5).Range.Text = "41÷3=13, 2"

$t.Cell(5, 1).Range.Text = "66÷6=11, 0"
$t.Cell(5, 2).Range.Text = "34÷2=17, 0"
$t.Cell(5, 3).Range.Text = "44÷3=14, 2"
$t.Cell(5, 4).Range.Text = "52÷2=26, 0"
$t.Cell(5, 5).Range.Text = "16÷9=1, 7"

$t.Cell(9, 1).Range.Text = "46÷9=5, 1"
$t.Cell(9, 2).Range.Text = "86÷4=21, 2"
$t.Cell(9, 3).Range.Text = "44÷9=4, 8"
$t.Cell(9, 4).Range.Text = "62÷7=8, 6"
$t.Cell(9, 5).Range.Text = "60÷2=30, 0"

$t.Cell(13, 1).Range.Text = "48÷6=8, 0"
$t.Cell(13, 2).Range.Text = "94÷2=47, 0"
$t.Cell(13, 3).Range.Text = "26÷7=3, 5"
$t.Cell(13, 4).Range.Text = "98÷3=32, 2"
$t.Cell(13, 5).Range.Text = "89÷2=44, 1"

$t.Cell(17, 1).Range.Text = "40÷3=13, 1"
$t.Cell(17, 2).Range.Text = "36÷2=18, 0"
$t.Cell(17, 3).Range.Text = "12÷6=2, 0"
$t.Cell(17, 4).Range.Text = "68÷5=13, 3"
$t.Cell(17, 5).Range.Text = "50÷8=6, 2"
